# Timings.xlsx: add Sec1/Sec2/Sec3/Average/Frac columns, a "Cache NE" row
# pair, and re-derive the Time column from the new Average/Frac chain
# instead of straight off the raw seconds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the old D column (formerly "Time", TEXT(...) formula) -----
# Its cells move to H and the formula chain changes shape, so wipe
# content+format first rather than leave stale style ids behind.
$ws.Range("D2:D5").Clear()

# --- header row -------------------------------------------------------
# (order matters for shared-string table layout: Time/Steps/Time Per
# Step get interned before the Sec1..Frac run, matching the source edit)
$ws.Range("H1").Value = "Time"
$ws.Range("I1").Value = "Steps"
$ws.Range("J1").Value = "Time Per Step"
$ws.Range("C1").Value = "Sec1"
$ws.Range("D1").Value = "Sec2"
$ws.Range("E1").Value = "Sec3"
$ws.Range("F1").Value = "Average"
$ws.Range("G1").Value = "Frac"

# --- row 2 : Debug / No Cache ------------------------------------------
$ws.Range("C2").Value = 179

# --- row 3 : Release / No Cache -----------------------------------------
$ws.Range("C3").Value = 144
$ws.Range("D3").Value = 146
$ws.Range("E3").Value = 136
$ws.Range("I3").Value = 1658103

# --- row 4 : Debug / Cache ----------------------------------------------
$ws.Range("C4").Value = 151

# --- row 5 : Release / Cache --------------------------------------------
$ws.Range("C5").Value = 101
$ws.Range("D5").Value = 99
$ws.Range("E5").Value = 104
$ws.Range("I5").Value = 1488535

# --- row 6 : Debug / Cache NE (new row) ----------------------------------
$ws.Range("A6").Value = "Debug"
$ws.Range("B6").Value = "Cache NE"
$ws.Range("C6").Value = 3
$ws.Range("I6").Value = 45277

# --- row 7 : Release / Cache NE (new row) --------------------------------
$ws.Range("A7").Value = "Release"
$ws.Range("B7").Value = "Cache NE"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("I7").Value = 45277

# --- formulas (written in the same order Excel itself grouped them so
# the shared-formula si indices line up: F si=0, G si=1, H si=2, F si=3)
$ws.Range("F2:F4").Formula = "=AVERAGE(C2:E2)"
$ws.Range("G2").Formula = "=F2/86400"
$ws.Range("H2").Formula = "=TEXT(G2,""hh:mm:ss"")"

$ws.Range("G3:G7").Formula = "=F3/86400"
$ws.Range("H3:H7").Formula = "=TEXT(G3,""hh:mm:ss"")"

$ws.Range("F5").Formula = "=AVERAGE(C5:E5)"
$ws.Range("J5").Formula = "=C5/I5"

$ws.Range("F6:F7").Formula = "=AVERAGE(C6:E6)"

# J6/J7/G8 stay blank but carry the fractional number format below.

# --- number formats -----------------------------------------------------
# Frac column (and the Steps ratio in J) gets a 7dp custom format; the
# averaged-seconds column F gets a 1dp custom format. H keeps reusing the
# workbook's pre-existing "hh:mm:ss AM/PM" custom format (style id 1) even
# though its formula now returns text via TEXT(...).
$ws.Range("G2:G8").NumberFormat = "0.0000000"
$ws.Range("J5:J7").NumberFormat = "0.0000000"
$ws.Range("F2:F7").NumberFormat = "0.0"
$ws.Range("H2:H7").NumberFormat = '[$-F400]h:mm:ss\ AM/PM'

# --- column widths --------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7
$ws.Columns.Item(2).ColumnWidth = 8.5703125
$ws.Range("C1:E1").ColumnWidth = 4.76171875
$ws.Columns.Item(6).ColumnWidth = 12.1875
$ws.Columns.Item(7).ColumnWidth = 9.1875
$ws.Columns.Item(8).ColumnWidth = 11.85546875
$ws.Columns.Item(10).ColumnWidth = 10.1875

# --- selection ------------------------------------------------------------
$ws.Range("G9").Select()
